# Scheduled market-data / leve-profit refresh
# Updates currentAveragePrice(NQ/HQ) and derived Leve price / profit
# columns (H:N) across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
# to reflect the latest Universalis market snapshot.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 247.76471
$ws.Range("I2").Value = 236.75
$ws.Range("K2").Value = 236.75
$ws.Range("M2").Value = -123.75

$ws.Range("H53").Value = 719.3333
$ws.Range("I53").Value = 729.7273
$ws.Range("J53").Value = 703
$ws.Range("K53").Value = 729.7273
$ws.Range("L53").Value = 703
$ws.Range("M53").Value = -92.72730000000001
$ws.Range("N53").Value = -1977

$ws.Range("H62").Value = 4774.8335
$ws.Range("I62").Value = 3051.25
$ws.Range("K62").Value = 3051.25
$ws.Range("M62").Value = -2427.25

$ws.Range("H65").Value = 4774.8335
$ws.Range("I65").Value = 3051.25
$ws.Range("K65").Value = 15256.25
$ws.Range("M65").Value = -12136.25

$ws.Range("H86").Value = 3262.8572
$ws.Range("I86").Value = 2560
$ws.Range("K86").Value = 2560
$ws.Range("M86").Value = -1437

$ws.Range("H89").Value = 3262.8572
$ws.Range("I89").Value = 2560
$ws.Range("K89").Value = 12800
$ws.Range("M89").Value = -7184

$ws.Range("H132").Value = 1832.6666
$ws.Range("I132").Value = 1385.2142
$ws.Range("K132").Value = 4155.642599999999
$ws.Range("M132").Value = -1625.642599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3343.1177
$ws.Range("I32").Value = 3343.1177
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3343.1177
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -3056.1177
$ws.Range("N32").ClearContents()

$ws.Range("H76").Value = 54000
$ws.Range("J76").Value = 54000
$ws.Range("L76").Value = 54000
$ws.Range("N76").Value = -54676

$ws.Range("H79").Value = 54000
$ws.Range("J79").Value = 54000
$ws.Range("L79").Value = 54000
$ws.Range("N79").Value = -56340

$ws.Range("H88").Value = 2768.8333
$ws.Range("I88").Value = 2396.125
$ws.Range("K88").Value = 2396.125
$ws.Range("M88").Value = -1990.125

$ws.Range("H91").Value = 2768.8333
$ws.Range("I91").Value = 2396.125
$ws.Range("K91").Value = 2396.125
$ws.Range("M91").Value = -992.125

$ws.Range("H110").Value = 4406.091
$ws.Range("I110").Value = 4995.875
$ws.Range("J110").Value = 2833.3333
$ws.Range("K110").Value = 4995.875
$ws.Range("L110").Value = 2833.3333
$ws.Range("M110").Value = -2950.875
$ws.Range("N110").Value = -6923.3333

$ws.Range("H132").Value = 916.65
$ws.Range("I132").Value = 916.65
$ws.Range("K132").Value = 2749.95
$ws.Range("M132").Value = -219.9499999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1982.5
$ws.Range("J20").Value = 2081
$ws.Range("L20").Value = 2081
$ws.Range("N20").Value = -2575

$ws.Range("H86").Value = 2145.2856
$ws.Range("I86").Value = 2602.2
$ws.Range("J86").Value = 1003
$ws.Range("K86").Value = 2602.2
$ws.Range("L86").Value = 1003
$ws.Range("M86").Value = -1479.2
$ws.Range("N86").Value = -3249

$ws.Range("H89").Value = 2145.2856
$ws.Range("I89").Value = 2602.2
$ws.Range("J89").Value = 1003
$ws.Range("K89").Value = 13011
$ws.Range("L89").Value = 5015
$ws.Range("M89").Value = -7395
$ws.Range("N89").Value = -16247

$ws.Range("H105").Value = 4683.1113
$ws.Range("I105").Value = 4235.4287
$ws.Range("K105").Value = 4235.4287
$ws.Range("M105").Value = -2488.4287

$ws.Range("H107").Value = 626.9
$ws.Range("I107").Value = 554.8333
$ws.Range("K107").Value = 554.8333
$ws.Range("M107").Value = 1365.1667

$ws.Range("H134").Value = 5951
$ws.Range("I134").Value = 5933.1816
$ws.Range("K134").Value = 17799.5448
$ws.Range("M134").Value = -15264.5448

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 2346.2632
$ws.Range("I7").Value = 1470.6666
$ws.Range("J7").Value = 5629.75
$ws.Range("K7").Value = 1470.6666
$ws.Range("L7").Value = 5629.75
$ws.Range("M7").Value = -1357.6666
$ws.Range("N7").Value = -5855.75

$ws.Range("H31").Value = 2572.7856
$ws.Range("J31").Value = 8999
$ws.Range("L31").Value = 8999
$ws.Range("N31").Value = -9589

$ws.Range("H34").Value = 2572.7856
$ws.Range("J34").Value = 8999
$ws.Range("L34").Value = 8999
$ws.Range("N34").Value = -9403

$ws.Range("H99").Value = 4669.4614
$ws.Range("J99").Value = 5567.3335
$ws.Range("L99").Value = 5567.3335
$ws.Range("N99").Value = -8563.333500000001

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H126").Value = 4669.4614
$ws.Range("J126").Value = 5567.3335
$ws.Range("L126").Value = 16702.0005
$ws.Range("N126").Value = -21642.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 635.8182
$ws.Range("I23").Value = 413
$ws.Range("K23").Value = 1239
$ws.Range("M23").Value = -1004

$ws.Range("H134").Value = 997.6667
$ws.Range("I134").Value = 997.6667
$ws.Range("K134").Value = 2993.0001
$ws.Range("M134").Value = 2076.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 696.125
$ws.Range("I2").Value = 952.9375
$ws.Range("J2").Value = 182.5
$ws.Range("K2").Value = 952.9375
$ws.Range("L2").Value = 182.5
$ws.Range("M2").Value = -839.9375
$ws.Range("N2").Value = -408.5

$ws.Range("H107").Value = 724.1429000000001
$ws.Range("I107").Value = 511.5
$ws.Range("K107").Value = 511.5
$ws.Range("M107").Value = 1408.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2244.389
$ws.Range("I7").Value = 1916.9231
$ws.Range("K7").Value = 1916.9231
$ws.Range("M7").Value = -1804.9231

$ws.Range("H22").Value = 2164.0908
$ws.Range("I22").Value = 1641
$ws.Range("J22").Value = 2600
$ws.Range("K22").Value = 1641
$ws.Range("L22").Value = 2600
$ws.Range("M22").Value = -1346
$ws.Range("N22").Value = -3190

$ws.Range("H27").Value = 2164.0908
$ws.Range("I27").Value = 1641
$ws.Range("J27").Value = 2600
$ws.Range("K27").Value = 1641
$ws.Range("L27").Value = 2600
$ws.Range("M27").Value = -1534
$ws.Range("N27").Value = -2814

$ws.Range("H40").Value = 1505.7273
$ws.Range("I40").Value = 1505.7273
$ws.Range("K40").Value = 1505.7273
$ws.Range("M40").Value = -1369.7273

$ws.Range("H43").Value = 12500
$ws.Range("I43").Value = 10000
$ws.Range("J43").Value = 15000
$ws.Range("K43").Value = 10000
$ws.Range("L43").Value = 15000
$ws.Range("M43").Value = -9807
$ws.Range("N43").Value = -15386

$ws.Range("H61").Value = 2330.6365
$ws.Range("I61").Value = 2106.8572
$ws.Range("J61").Value = 2722.25
$ws.Range("K61").Value = 2106.8572
$ws.Range("L61").Value = 2722.25
$ws.Range("M61").Value = -1904.8572
$ws.Range("N61").Value = -3126.25

$ws.Range("H113").Value = 2330.6365
$ws.Range("I113").Value = 2106.8572
$ws.Range("J113").Value = 2722.25
$ws.Range("K113").Value = 2106.8572
$ws.Range("L113").Value = 2722.25
$ws.Range("M113").Value = 63.14280000000008
$ws.Range("N113").Value = -7062.25

$ws.Range("H122").Value = 5496.75
$ws.Range("I122").Value = 5458.6924
$ws.Range("J122").Value = 5661.6665
$ws.Range("K122").Value = 16376.0772
$ws.Range("L122").Value = 16984.9995
$ws.Range("M122").Value = -13926.0772
$ws.Range("N122").Value = -21884.9995

$ws.Range("H126").Value = 2244.389
$ws.Range("I126").Value = 1916.9231
$ws.Range("K126").Value = 5750.7693
$ws.Range("M126").Value = -3280.7693

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2107.6365
$ws.Range("I132").Value = 1873
$ws.Range("K132").Value = 5619
$ws.Range("M132").Value = -3089

$ws.Range("H136").Value = 1744.1875
$ws.Range("J136").Value = 1099.5
$ws.Range("L136").Value = 3298.5
$ws.Range("N136").Value = -8398.5
